$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the lecture room for 12.09 (cell F4: Seminar column) from BORCH to AUD N
$ws.Range("F4").Value = "12.09 *08:15 - 10:00* (**AUD N**)"

# Move the active selection to F5 (reflects the saved view state in the diff)
$ws.Range("F5").Select()
